{"js": "// Replace the date and the multiplication-table entries with their\n// updated values. Every source string in this document is unique, so a\n// straightforward search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2025-01-25 Saturday\", \"2025-01-26 Sunday\"],\n  [\"12\u00d747=564\", \"32\u00d781=2592\"],\n  [\"96\u00d793=8928\", \"43\u00d747=2021\"],\n  [\"76\u00d729=2204\", \"13\u00d711=143\"],\n  [\"26\u00d783=2158\", \"52\u00d712=624\"],\n  [\"64\u00d753=3392\", \"62\u00d785=5270\"],\n  [\"52\u00d778=4056\", \"50\u00d783=4150\"],\n  [\"42\u00d721=882\", \"67\u00d758=3886\"],\n  [\"21\u00d738=798\", \"65\u00d769=4485\"],\n  [\"85\u00d760=5100\", \"91\u00d733=3003\"],\n  [\"27\u00d791=2457\", \"80\u00d788=7040\"],\n  [\"64\u00d773=4672\", \"77\u00d739=3003\"],\n  [\"63\u00d795=5985\", \"25\u00d714=350\"],\n  [\"68\u00d725=1700\", \"98\u00d742=4116\"],\n  [\"52\u00d737=1924\", \"76\u00d732=2432\"],\n  [\"26\u00d711=286\", \"13\u00d766=858\"],\n  [\"26\u00d752=1352\", \"55\u00d783=4565\"],\n  [\"98\u00d785=8330\", \"92\u00d787=8004\"],\n  [\"34\u00d740=1360\", \"98\u00d750=4900\"],\n  [\"13\u00d723=299\", \"54\u00d763=3402\"],\n  [\"82\u00d794=7708\", \"13\u00d721=273\"],\n  [\"31\u00d721=651\", \"70\u00d791=6370\"],\n  [\"92\u00d779=7268\", \"44\u00d771=3124\"],\n  [\"49\u00d740=1960\", \"28\u00d716=448\"],\n  [\"81\u00d780=6480\", \"46\u00d776=3496\"],\n  [\"94\u00d722=2068\", \"68\u00d730=2040\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and every multiplication fact in the practice\n# table. Each \"old\" string below occurs exactly once in the document, so\n# a plain Find/Replace (wdReplaceAll) per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-25 Saturday\", \"2025-01-26 Sunday\"),\n    @(\"12\u00d747=564\", \"32\u00d781=2592\"),\n    @(\"96\u00d793=8928\", \"43\u00d747=2021\"),\n    @(\"76\u00d729=2204\", \"13\u00d711=143\"),\n    @(\"26\u00d783=2158\", \"52\u00d712=624\"),\n    @(\"64\u00d753=3392\", \"62\u00d785=5270\"),\n    @(\"52\u00d778=4056\", \"50\u00d783=4150\"),\n    @(\"42\u00d721=882\", \"67\u00d758=3886\"),\n    @(\"21\u00d738=798\", \"65\u00d769=4485\"),\n    @(\"85\u00d760=5100\", \"91\u00d733=3003\"),\n    @(\"27\u00d791=2457\", \"80\u00d788=7040\"),\n    @(\"64\u00d773=4672\", \"77\u00d739=3003\"),\n    @(\"63\u00d795=5985\", \"25\u00d714=350\"),\n    @(\"68\u00d725=1700\", \"98\u00d742=4116\"),\n    @(\"52\u00d737=1924\", \"76\u00d732=2432\"),\n    @(\"26\u00d711=286\", \"13\u00d766=858\"),\n    @(\"26\u00d752=1352\", \"55\u00d783=4565\"),\n    @(\"98\u00d785=8330\", \"92\u00d787=8004\"),\n    @(\"34\u00d740=1360\", \"98\u00d750=4900\"),\n    @(\"13\u00d723=299\", \"54\u00d763=3402\"),\n    @(\"82\u00d794=7708\", \"13\u00d721=273\"),\n    @(\"31\u00d721=651\", \"70\u00d791=6370\"),\n    @(\"92\u00d779=7268\", \"44\u00d771=3124\"),\n    @(\"49\u00d740=1960\", \"28\u00d716=448\"),\n    @(\"81\u00d780=6480\", \"46\u00d776=3496\"),\n    @(\"94\u00d722=2068\", \"68\u00d730=2040\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
